$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: header "Save" in H1 (reuse the same formatting as the other
# header cells, e.g. G1, via copy/paste-format), and 0 in H2:H9.

$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
